# Update column G ("K" = strikeouts replaced by K count) values for rows 2-15
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 3
    14 = 1
    15 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
